# v3.0 update FCI 27/1/2023
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows, sorted alphabetically by fund name, with avg/total moved to the bottom.
# Column B keeps the existing (06-01-2023) values; column C adds the new (13-01-2023) values.
$rows = @(
    @{ Label = "1810 Renta variable";       B = 49975.45;             C = 49943.79 },
    @{ Label = "Alpha Acciones";            B = 92655.67999999999;    C = 92648.34 },
    @{ Label = "Alpha Mega";                B = 50825.22;             C = 50658.61 },
    @{ Label = "Delta Acciones";            B = 22975.88;             C = 22951.24 },
    @{ Label = "Delta Recursos Naturales";  B = 63076.57;             C = 62902.29 },
    @{ Label = "Delta Select";              B = 42184.27;             C = 42202.12 },
    @{ Label = "Delta gestion V";           B = 9677.809999999999;    C = 9968.32 },
    @{ Label = "Fima Acciones";             B = 181878.46;            C = 195576.51 },
    @{ Label = "Fima PB Acciones";          B = 96412.03;             C = 101172.07 },
    @{ Label = "Gainvest Renta Variable";   B = 38574.04;             C = 47376.57 },
    @{ Label = "HF Acciones Argentinas";    B = 901.3099999999999;    C = 943.59 },
    @{ Label = "avg";                       B = 59012.43;             C = 61485.77 },
    @{ Label = "total";                     B = 649136.72;            C = 676343.45 }
)

# Header row: add the new date in column C (column B header already holds 06-01-2023)
$ws.Cells.Item(1, 3).Value = "13-01-2023"
$ws.Range("B1").Copy()
$ws.Range("C1").PasteSpecial(-4122)  # xlPasteFormats, matches B1's bold/border/alignment style

# Write the data rows starting at row 2, rewriting labels/values in the new order
$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row.Label
    $ws.Cells.Item($r, 2).Value = $row.B
    $ws.Cells.Item($r, 3).Value = $row.C
    $r = $r + 1
}

# Column A keeps its existing bold/border style in every row (already style-1 in the source
# sheet); make sure it stays applied after the label rewrite above.
$ws.Range("A2").Copy()
$ws.Range("A3:A14").PasteSpecial(-4122)  # xlPasteFormats
